# Auto-generated edit script applying cell text updates per diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '26.740.84'
Set-TextValue $ws.Range('E2') '  -0.70%  '
Set-TextValue $ws.Range('D3') '1.546.33'
Set-TextValue $ws.Range('E3') '  -1.26%  '
Set-TextValue $ws.Range('E4') '  -0.28%  '
Set-TextValue $ws.Range('D5') '206.37'
Set-TextValue $ws.Range('E5') '  +0.08%  '
Set-TextValue $ws.Range('E6') '  -1.38%  '
Set-TextValue $ws.Range('E7') '  -0.15%  '
Set-TextValue $ws.Range('D8') '21.41'
Set-TextValue $ws.Range('E8') '  -3.25%  '
Set-TextValue $ws.Range('E9') '  -1.17%  '
Set-TextValue $ws.Range('E10') '  -0.57%  '
Set-TextValue $ws.Range('D11') '0.0854'
Set-TextValue $ws.Range('E11') '  -1.45%  '
Set-TextValue $ws.Range('E12') '  -1.31%  '
Set-TextValue $ws.Range('D13') '1.546.41'
Set-TextValue $ws.Range('E13') '  -1.37%  '
Set-TextValue $ws.Range('D14') '3.68'
Set-TextValue $ws.Range('E14') '  -2.11%  '
Set-TextValue $ws.Range('E15') '  -1.02%  '
Set-TextValue $ws.Range('D16') '26.699.73'
Set-TextValue $ws.Range('E16') '  -0.90%  '
Set-TextValue $ws.Range('D17') '61.20'
Set-TextValue $ws.Range('E17') '  -0.95%  '
Set-TextValue $ws.Range('D18') '212.50'
Set-TextValue $ws.Range('E18') '  -0.93%  '
Set-TextValue $ws.Range('D19') '0.0₃0689'
Set-TextValue $ws.Range('E19') '  +1.28%  '
Set-TextValue $ws.Range('D20') '7.24'
Set-TextValue $ws.Range('E20') '  -1.56%  '
Set-TextValue $ws.Range('E21') '  -0.25%  '
Set-TextValue $ws.Range('E22') '  -1.03%  '
Set-TextValue $ws.Range('D23') '8.98'
Set-TextValue $ws.Range('E23') '  -4.12%  '
Set-TextValue $ws.Range('E24') '  -1.00%  '
Set-TextValue $ws.Range('D25') '152.49'
Set-TextValue $ws.Range('E25') '  +0.42%  '
Set-TextValue $ws.Range('D26') '14.87'
Set-TextValue $ws.Range('E26') '  +0.05%  '
Set-TextValue $ws.Range('D27') '6.47'
Set-TextValue $ws.Range('E27') '  -3.12%  '
Set-TextValue $ws.Range('E28') '  -0.25%  '
Set-TextValue $ws.Range('E29') '  -0.90%  '
Set-TextValue $ws.Range('B30') 'PancakeSwap'
Set-TextValue $ws.Range('C30') 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue $ws.Range('D30') '1.10'
Set-TextValue $ws.Range('E30') '  -0.96%  '
Set-TextValue $ws.Range('B31') 'Hedera'
Set-TextValue $ws.Range('C31') 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue $ws.Range('D31') '0.0458'
Set-TextValue $ws.Range('E31') '  -0.50%  '
Set-TextValue $ws.Range('E32') '  +1.51%  '
Set-TextValue $ws.Range('D33') '1.351.66'
Set-TextValue $ws.Range('E33') '  -2.49%  '
Set-TextValue $ws.Range('D34') '2.92'
Set-TextValue $ws.Range('E34') '  +0.51%  '
Set-TextValue $ws.Range('E35') '  -3.00%  '
Set-TextValue $ws.Range('E36') '  -0.85%  '
Set-TextValue $ws.Range('D37') '0.935'
Set-TextValue $ws.Range('E37') '  -0.50%  '
Set-TextValue $ws.Range('E38') '  +0.33%  '
Set-TextValue $ws.Range('E39') '  +2.35%  '
Set-TextValue $ws.Range('D40') '0.801'
Set-TextValue $ws.Range('E40') '  -1.06%  '
Set-TextValue $ws.Range('E41') '  +5.23%  '
Set-TextValue $ws.Range('D42') '0.995'
Set-TextValue $ws.Range('E42') '  -0.73%  '
Set-TextValue $ws.Range('E43') '  +0.16%  '
Set-TextValue $ws.Range('E44') '  -2.84%  '
Set-TextValue $ws.Range('D45') '62.60'
Set-TextValue $ws.Range('E45') '  -1.33%  '
Set-TextValue $ws.Range('D46') '1.678.14'
Set-TextValue $ws.Range('E46') '  -1.30%  '
Set-TextValue $ws.Range('E47') '  -4.57%  '
Set-TextValue $ws.Range('D48') '85.89'
Set-TextValue $ws.Range('E48') '  +0.64%  '
Set-TextValue $ws.Range('D49') '0.0509'
Set-TextValue $ws.Range('E49') '  +2.94%  '
Set-TextValue $ws.Range('B50') 'Algorand'
Set-TextValue $ws.Range('C50') 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue $ws.Range('D50') '0.0952'
Set-TextValue $ws.Range('E50') '  +0.58%  '
Set-TextValue $ws.Range('B51') 'BabyDogeCoin'
Set-TextValue $ws.Range('C51') 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextValue $ws.Range('D51') '0.0₇0953'
Set-TextValue $ws.Range('E51') '  -1.76%  '
